$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new column before D ("P (W)" - power column)
# ---------------------------------------------------------------
$ws.Columns("D:D").Insert()
$ws.Columns("D:D").ColumnWidth = 6

# Header for the new column
$ws.Range("D1").Value = "P (W)"

# Rename old "Output" header (now in column F after the shift) to "I/O"
$ws.Range("F1").Value = "I/O"

# ---------------------------------------------------------------
# 2. Insert a new row before row 9 (Front_Alive) for the new
#    "Tyre_temperatuae" sensor, and copy formatting from the row
#    above so it matches the rest of the table.
# ---------------------------------------------------------------
$ws.Rows("9:9").Insert()
$ws.Range("A8:H8").Copy()
$ws.Range("A9:H9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A9").Value = "Tyre_temperatuae"
$ws.Range("B9").Value = "MLX90641"
$ws.Range("C9").Value = 3.3
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = "I2C"
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = "https://drive.google.com/file/d/1QEUmVspiburSFoTysupAeOcM1dWjx5MP/view"

# New "P (W)" value cell gets a numeric 0.000 format (value left blank)
$ws.Range("D9").NumberFormat = "0.000"

# ---------------------------------------------------------------
# 3. Decorative border box around H8:I11 (left border column I,
#    left+top+bottom border column H) framing the new sensor row.
# ---------------------------------------------------------------
foreach ($r in 8..11) {
    $ws.Cells.Item($r, 8).Borders.Item(7).LineStyle = 1
    $ws.Cells.Item($r, 8).Borders.Item(8).LineStyle = 1
    $ws.Cells.Item($r, 8).Borders.Item(9).LineStyle = 1
    $ws.Cells.Item($r, 8).Borders.Item(10).LineStyle = -4142
    $ws.Cells.Item($r, 9).Borders.Item(7).LineStyle = 1
    $ws.Cells.Item($r, 9).Borders.Item(8).LineStyle = -4142
    $ws.Cells.Item($r, 9).Borders.Item(9).LineStyle = -4142
    $ws.Cells.Item($r, 9).Borders.Item(10).LineStyle = -4142
}

# ---------------------------------------------------------------
# 4. Final selection / view state
# ---------------------------------------------------------------
$ws.Range("D9").Select()
